$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet and
#    populate it with the per-fund holding detail for 2022-Q1, mirroring the
#    layout used by the other quarterly sheets (2020-Q4 .. 2021-Q4).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Make header + index column look like the other quarterly sheets: bold font,
# thin border all around, centered/top aligned.
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$idxRange = $q1.Range("A2:A4")
$idxRange.Font.Bold = $true
$idxRange.Borders.LineStyle = 1
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160

# Columns B, D, E, F, G hold numeric-looking values that must stay TEXT
# (leading-zero fund codes, fixed decimal strings) - force text format before
# writing so Excel doesn't silently coerce them to numbers.
$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

# Row 2 - 006199
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "006199"
$q1.Range("C2").Value = "长盛同锦研究精选混合"
$q1.Range("D2").Value = "1.73"
$q1.Range("E2").Value = "82.48"
$q1.Range("F2").Value = "2.24"
$q1.Range("G2").Value = "0.0388"
$q1.Range("H2").Value = 10

# Row 3 - 001892
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "001892"
$q1.Range("C3").Value = "长盛新兴成长主题灵活配置混合"
$q1.Range("D3").Value = "1.32"
$q1.Range("E3").Value = "82.10"
$q1.Range("F3").Value = "2.24"
$q1.Range("G3").Value = "0.0296"
$q1.Range("H3").Value = 10

# Row 4 - 002085
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "002085"
$q1.Range("C4").Value = "长盛互联网+主题灵活配置混合"
$q1.Range("D4").Value = "0.84"
$q1.Range("E4").Value = "83.97"
$q1.Range("F4").Value = "2.24"
$q1.Range("G4").Value = "0.0188"
$q1.Range("H4").Value = 10

# Match the page margins used by every other sheet in the workbook (values
# are expressed in points: 0.75in=54, 1in=72, 0.5in=36).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2) Prepend a 2022-Q1 summary row to the "总计" sheet, shifting the existing
#    rows down by one (and keeping column A's running index in B2..B? order).
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows("2:2").Insert()
$totals.Range("A2:D2").ClearFormats()

$totals.Range("A2").Font.Bold = $true
$totals.Range("A2").Borders.LineStyle = 1
$totals.Range("A2").HorizontalAlignment = -4108
$totals.Range("A2").VerticalAlignment = -4160

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 0.09
